$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62; existing rows 62:159 shift down to 63:160.
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with its data. Columns that are
# constant across the whole table (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) are filled
# with the same values as the surrounding rows; D,M,N,O,P,S carry the new
# record's values.
$ws.Range("A62").Value = 9
$ws.Range("B62").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C62").Value = "Metropolitana"
$ws.Range("D62").Value = 44580
$ws.Range("E62").Value = 13
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100101
$ws.Range("H62").Value = "Berries"
$ws.Range("I62").Value = 100101001
$ws.Range("J62").Value = "Arándano (blue)"
$ws.Range("K62").Value = "Sin especificar"
$ws.Range("L62").Value = "Primera"
$ws.Range("M62").Value = 410
$ws.Range("N62").Value = 4000
$ws.Range("O62").Value = 4000
$ws.Range("P62").Value = 4000
$ws.Range("Q62").Value = "`$/bandeja 2 kilos"
$ws.Range("R62").Value = "Región de O'Higgins"
$ws.Range("S62").Value = 2000
$ws.Range("T62").Value = 2
